# "adding team gradebook as dp2 makeup"
# On slide 1, the "Today's Attendance password" textbox has a second
# paragraph that is just a blank-filled line ("__________") waiting for
# the day's password to be written in. Replace it with this session's
# password, "morefun!".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox shape (id=4, name "TextBox 3") that holds the
# attendance password instead of assuming a fixed Shapes index.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 4) {
        $shp = $candidate
    }
}

$tr = $shp.TextFrame.TextRange

# Find the paragraph holding the blank-line placeholder.
$pwdPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidatePara = $tr.Paragraphs($i)
    if ($candidatePara.Text -like "*____*") {
        $pwdPara = $candidatePara
    }
}

$blankRun = $pwdPara.Runs(1)

# Split the blank run into "morefun" + "!" (two runs), preserving the
# sz=4400 / yellow-highlight formatting already on it.
$blankRun.InsertAfter("!") | Out-Null
$firstRun = $pwdPara.Runs(1)
$firstRun.Text = "morefun"
